$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.492.33"
$ws.Range("E2").Value = "  +2.01%  "

# Row 3
$ws.Range("D3").Value = "1.990.89"
$ws.Range("E3").Value = "  +5.92%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").Value = "'325.38"
$ws.Range("E5").Value = "  +0.15%  "

# Row 6
$ws.Range("D6").Value = "'0.9993"
$ws.Range("E6").Value = "  -0.25%  "

# Row 7
$ws.Range("D7").Value = "'0.4674"
$ws.Range("E7").Value = "  +1.53%  "

# Row 8
$ws.Range("D8").Value = "'0.3945"
$ws.Range("E8").Value = "  +1.75%  "

# Row 9
$ws.Range("D9").Value = "'46.39"
$ws.Range("E9").Value = "  -0.51%  "

# Row 10
$ws.Range("D10").Value = "'0.07933"
$ws.Range("E10").Value = "  +0.94%  "

# Row 11
$ws.Range("D11").Value = "'1.000"
$ws.Range("E11").Value = "  +1.42%  "

# Row 12
$ws.Range("D12").Value = "'22.87"
$ws.Range("E12").Value = "  +5.05%  "

# Row 13
$ws.Range("D13").Value = "1.972.92"
$ws.Range("E13").Value = "  +11.31%  "

# Row 14
$ws.Range("D14").Value = "'7.272"
$ws.Range("E14").Value = "  +3.89%  "

# Row 15
$ws.Range("D15").Value = "'5.861"
$ws.Range("E15").Value = "  +3.76%  "

# Row 16
$ws.Range("D16").Value = "'0.07125"
$ws.Range("E16").Value = "  +2.50%  "

# Row 17
$ws.Range("D17").Value = "'88.58"
$ws.Range("E17").Value = "  +0.55%  "

# Row 18
$ws.Range("E18").Value = "  -0.16%  "

# Row 19
$ws.Range("D19").Value = "'0.000009963"
$ws.Range("E19").Value = "  -0.14%  "

# Row 20
$ws.Range("E20").Value = "  +2.23%  "

# Row 21
$ws.Range("D21").Value = "'0.9993"
$ws.Range("E21").Value = "  -0.18%  "

# Row 22
$ws.Range("D22").Value = "29.595.84"
$ws.Range("E22").Value = "  +2.43%  "

# Row 23
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'5.530"
$ws.Range("E23").Value = "  +5.68%  "

# Row 24
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "'11.26"
$ws.Range("E24").Value = "  +2.71%  "

# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'2.101"
$ws.Range("E25").Value = "  +0.71%  "

# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'157.62"
$ws.Range("E26").Value = "  +0.79%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'19.63"
$ws.Range("E27").Value = "  +1.76%  "

# Row 28
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "'5.973"
$ws.Range("E28").Value = "  -0.44%  "

# Row 29
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "'120.19"
$ws.Range("E29").Value = "  +2.34%  "

# Row 30
$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "'1.956"
$ws.Range("E30").Value = "  +1.48%  "

# Row 31
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.09441"
$ws.Range("E31").Value = "  +0.79%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'0.9041"
$ws.Range("E32").Value = "  +0.09%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.251"
$ws.Range("E33").Value = "  -0.19%  "

# Row 34
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.347"
$ws.Range("E34").Value = "  +2.56%  "

# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'3.176"
$ws.Range("E35").Value = "  -2.40%  "

# Row 36
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "'0.000003538"
$ws.Range("E36").Value = "  +120.70%  "

# Row 37
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.05827"
$ws.Range("E37").Value = "  +1.45%  "

# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.173"
$ws.Range("E38").Value = "  -0.70%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02114"
$ws.Range("E39").Value = "  +2.09%  "

# Row 40
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'7.854"
$ws.Range("E40").Value = "  +3.00%  "

# Row 41
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5746"
$ws.Range("E41").Value = "  +1.70%  "

# Row 42
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.1824"
$ws.Range("E42").Value = "  +3.42%  "

# Row 43
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'9.796"
$ws.Range("E43").Value = "  +1.25%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'12.04"
$ws.Range("E44").Value = "  +1.10%  "

# Row 45
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.5368"
$ws.Range("E45").Value = "  +0.44%  "

# Row 46
$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").Value = "'2.687"
$ws.Range("E46").Value = "  +6.23%  "

# Row 47
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'2.158"
$ws.Range("E47").Value = "  -5.01%  "

# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.867"
$ws.Range("E48").Value = "  +1.27%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.06940"
$ws.Range("E49").Value = "  -1.48%  "

# Row 50
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'113.96"
$ws.Range("E50").Value = "  +1.13%  "

# Row 51
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").Value = "'0.3099"
$ws.Range("E51").Value = "  +8.55%  "

